# Fruta / hortaliza, semanal
#
# The weekly refresh re-sorts the price records (rows 2-24) of the single
# data sheet. Every record's full set of fields - Fecha (D), Variedad (K),
# Calidad (L), Volumen (M), Precio minimo/maximo/promedio (N/O/P), Unidad
# de comercializacion (Q), Origen (R), Precio $/Kg (S) and Kg/unidad (T) -
# moves together as one unit to a (possibly) different row; the leading
# descriptive columns (A-C, E-J: Mercado/Region/Codreg/Tipo/Producto/...)
# stay put since they describe the same market/product for every row.
#
# Row -> source-row mapping (new row N gets the record that used to live
# in row $mapping[N]):
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 10
    3  = 11
    4  = 3
    5  = 19
    6  = 20
    7  = 6
    8  = 7
    9  = 21
    10 = 2
    11 = 22
    12 = 23
    13 = 24
    14 = 12
    15 = 13
    16 = 16
    17 = 14
    18 = 17
    19 = 18
    20 = 8
    21 = 9
    22 = 15
    23 = 4
    24 = 5
}

$cols = @("D","K","L","M","N","O","P","Q","R","S","T")

# Snapshot every source cell first so writes never clobber a value that
# still needs to be read later (several rows swap with each other).
$snapshot = @{}
foreach ($r in 2..24) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowData
}

foreach ($r in 2..24) {
    $srcRow = $mapping[$r]
    $rowData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $rowData[$c]
    }
}
